# Daily attendance processing - 2025-10-05 05:39:19
# Applies the attendance-report refresh: updated class-statistics counters,
# updated group-statistics row for Year 2 / A2, re-ordered "Recorded By"
# grader lists (same people, different order) across many rows, and the
# PARASITOLOGY / A2 / Session 2 row (row 25) moving from "Pending" to
# "Recorded" (value + formatting + status).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a percentage string into a cell as literal TEXT (matching
# the source file's t="inlineStr" percentage cells) without disturbing the
# destination cell's existing style. A plain `.Value = "26.8%"` assignment
# gets auto-coerced by Excel into a percent NUMBER (and stamps a new percent
# number-format style onto the cell), so instead we stage the text in a
# scratch cell explicitly formatted as Text, then copy/paste-special just
# the values into the real target (format of the target is left alone).
$scratch = $ws.Range("Z1")
function Set-TextPercent($rangeAddress, $text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($rangeAddress).PasteSpecial(-4163)   # xlPasteValues
}

# ---------------------------------------------------------------------
# Class Statistics block (K/L columns near the top of the sheet)
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 41        # Recorded Sessions

$ws.Range("L8").Value = 104       # Pending Sessions

Set-TextPercent "L9" "26.8%"      # Coverage %

Set-TextPercent "L10" "48.9%"     # Average Attendance %

# ---------------------------------------------------------------------
# Group Statistics row for Year 2 / A2 (row 16)
# ---------------------------------------------------------------------
$ws.Range("O16").Value = 6        # Recorded

$ws.Range("Q16").Value = 10       # Pending

Set-TextPercent "R16" "35.3%"     # Coverage %

Set-TextPercent "S16" "57.1%"     # Avg Attendance %

$scratch.Clear()

# ---------------------------------------------------------------------
# Row 25: PARASITOLOGY / A2 / Session 2 — moved from Pending to Recorded.
# Copy the formatting from row 24 (an already-"Recorded" row) so the fill
# / font match the other recorded rows exactly, then set the new values.
# ---------------------------------------------------------------------
$ws.Range("A24:I24").Copy()
$ws.Range("A25:I25").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("G25").Value = "Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("H25").Value = "97/216"
$ws.Range("I25").Value = "Recorded"

# ---------------------------------------------------------------------
# "Recorded By" grader-list reorderings (same graders, new order) —
# these mirror how the upstream attendance system re-serializes the
# recorder list on each processing run.
# ---------------------------------------------------------------------
$ws.Range("G14").Value = "nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G17").Value = "ola.m.abdelfattah@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G31").Value = "nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G34").Value = "ola.m.abdelfattah@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G35").Value = "neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"
$ws.Range("G45").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, backup@backdoor.com, System"
$ws.Range("G51").Value = "yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"
$ws.Range("G62").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, backup@backdoor.com, System"
$ws.Range("G68").Value = "yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"
$ws.Range("G72").Value = "wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G76").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G83").Value = "afaf.abdallah@med.asu.edu.eg, Youstina.ibrahim@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G85").Value = "ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
$ws.Range("G98").Value = "afaf.abdallah@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg, user@user.com"
$ws.Range("G99").Value = "user@user.com, Walaa.h.ghanima@med.asu.edu.eg"
$ws.Range("G102").Value = "ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
$ws.Range("G116").Value = "afaf.abdallah@med.asu.edu.eg, enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G119").Value = "ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marinasorial@med.asu.edu.eg"
$ws.Range("G133").Value = "afaf.abdallah@med.asu.edu.eg, enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G136").Value = "ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marinasorial@med.asu.edu.eg"
$ws.Range("G149").Value = "user@user.com, Walaa.h.ghanima@med.asu.edu.eg"
$ws.Range("G150").Value = "afaf.abdallah@med.asu.edu.eg, Youstina.ibrahim@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
